$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 95.24255366666667
$ws.Range("H2").Value = 285.727661
$ws.Range("I2").Value = 0.2732032672746878
$ws.Range("J2").Value = 0.2732032672746877
$ws.Range("M2").Value = 52.91852733333334
$ws.Range("N2").Value = 158.755582
$ws.Range("O2").Value = 0.9912603569328422
$ws.Range("P2").Value = 0.9912603569328421
$ws.Range("Q2").Value = 5040.095679505967
$ws.Range("R2").Value = 45360.86111555371
$ws.Range("S2").Value = 0.2708155682339257
$ws.Range("T2").Value = 0.2708155682339256

$ws.Range("G3").Value = 95.24255366666667
$ws.Range("H3").Value = 285.727661
$ws.Range("I3").Value = 0.2732032672746878
$ws.Range("J3").Value = 0.2732032672746877
$ws.Range("O3").Value = 0.003851187374513192
$ws.Range("P3").Value = 0.003851187374513192
$ws.Range("Q3").Value = 19.581488063652
$ws.Range("R3").Value = 176.233392572868
$ws.Range("S3").Value = 0.001052156973604031
$ws.Range("T3").Value = 0.001052156973604031

$ws.Range("G4").Value = 95.24255366666667
$ws.Range("H4").Value = 285.727661
$ws.Range("I4").Value = 0.2732032672746878
$ws.Range("J4").Value = 0.2732032672746877
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2609706666666667
$ws.Range("N4").Value = 0.7829120000000001
$ws.Range("O4").Value = 0.004888455692644593
$ws.Range("P4").Value = 0.004888455692644592
$ws.Range("Q4").Value = 24.85551272542578
$ws.Range("R4").Value = 223.699614528832
$ws.Range("S4").Value = 0.00133554206715805
$ws.Range("T4").Value = 0.001335542067158049

$ws.Range("I5").Value = 0.3246226791565123
$ws.Range("J5").Value = 0.3246226791565122
$ws.Range("M5").Value = 52.91852733333334
$ws.Range("N5").Value = 158.755582
$ws.Range("O5").Value = 0.9912603569328422
$ws.Range("P5").Value = 0.9912603569328421
$ws.Range("Q5").Value = 5988.688857960728
$ws.Range("R5").Value = 53898.19972164654
$ws.Range("S5").Value = 0.3217855928091798
$ws.Range("T5").Value = 0.3217855928091798

$ws.Range("I6").Value = 0.3246226791565123
$ws.Range("J6").Value = 0.3246226791565122
$ws.Range("O6").Value = 0.003851187374513192
$ws.Range("P6").Value = 0.003851187374513192
$ws.Range("S6").Value = 0.001250182763448207
$ws.Range("T6").Value = 0.001250182763448207

$ws.Range("I7").Value = 0.3246226791565123
$ws.Range("J7").Value = 0.3246226791565122
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2609706666666667
$ws.Range("N7").Value = 0.7829120000000001
$ws.Range("O7").Value = 0.004888455692644593
$ws.Range("P7").Value = 0.004888455692644592
$ws.Range("Q7").Value = 29.53355284958578
$ws.Range("R7").Value = 265.801975646272
$ws.Range("S7").Value = 0.001586903583884192
$ws.Range("T7").Value = 0.001586903583884191

$ws.Range("G8").Value = 89.83461266666666
$ws.Range("H8").Value = 269.503838
$ws.Range("I8").Value = 0.2576905883979786
$ws.Range("J8").Value = 0.2576905883979785
$ws.Range("M8").Value = 52.91852733333334
$ws.Range("N8").Value = 158.755582
$ws.Range("O8").Value = 0.9912603569328422
$ws.Range("P8").Value = 0.9912603569328421
$ws.Range("Q8").Value = 4753.915405880412
$ws.Range("R8").Value = 42785.23865292371
$ws.Range("S8").Value = 0.2554384646336144
$ws.Range("T8").Value = 0.2554384646336143

$ws.Range("G9").Value = 89.83461266666666
$ws.Range("H9").Value = 269.503838
$ws.Range("I9").Value = 0.2576905883979786
$ws.Range("J9").Value = 0.2576905883979785
$ws.Range("O9").Value = 0.003851187374513192
$ws.Range("P9").Value = 0.003851187374513192
$ws.Range("Q9").Value = 18.469637025816
$ws.Range("R9").Value = 166.226733232344
$ws.Range("S9").Value = 0.0009924147405691709
$ws.Range("T9").Value = 0.0009924147405691705

$ws.Range("G10").Value = 89.83461266666666
$ws.Range("H10").Value = 269.503838
$ws.Range("I10").Value = 0.2576905883979786
$ws.Range("J10").Value = 0.2576905883979785
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2609706666666667
$ws.Range("N10").Value = 0.7829120000000001
$ws.Range("O10").Value = 0.004888455692644593
$ws.Range("P10").Value = 0.004888455692644592
$ws.Range("Q10").Value = 23.44419875736178
$ws.Range("R10").Value = 210.997788816256
$ws.Range("S10").Value = 0.001259709023795033
$ws.Range("T10").Value = 0.001259709023795032

$ws.Range("G11").Value = 50.36899566666667
$ws.Range("H11").Value = 151.106987
$ws.Range("I11").Value = 0.1444834651708214
$ws.Range("J11").Value = 0.1444834651708214
$ws.Range("M11").Value = 52.91852733333334
$ws.Range("N11").Value = 158.755582
$ws.Range("O11").Value = 0.9912603569328422
$ws.Range("P11").Value = 0.9912603569328421
$ws.Range("Q11").Value = 2665.453073939048
$ws.Range("R11").Value = 23989.07766545143
$ws.Range("S11").Value = 0.1432207312561223
$ws.Range("T11").Value = 0.1432207312561223

$ws.Range("G12").Value = 50.36899566666667
$ws.Range("H12").Value = 151.106987
$ws.Range("I12").Value = 0.1444834651708214
$ws.Range("J12").Value = 0.1444834651708214
$ws.Range("O12").Value = 0.003851187374513192
$ws.Range("P12").Value = 0.003851187374513192
$ws.Range("Q12").Value = 10.355664033084
$ws.Range("R12").Value = 93.200976297756
$ws.Range("S12").Value = 0.0005564328968917841
$ws.Range("T12").Value = 0.000556432896891784

$ws.Range("G13").Value = 50.36899566666667
$ws.Range("H13").Value = 151.106987
$ws.Range("I13").Value = 0.1444834651708214
$ws.Range("J13").Value = 0.1444834651708214
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2609706666666667
$ws.Range("N13").Value = 0.7829120000000001
$ws.Range("O13").Value = 0.004888455692644593
$ws.Range("P13").Value = 0.004888455692644592
$ws.Range("Q13").Value = 13.14483037846045
$ws.Range("R13").Value = 118.303473406144
$ws.Range("S13").Value = 0.0007063010178073188
$ws.Range("T13").Value = 0.0007063010178073186

